$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "63.032.63"
$ws.Range("E2").Value = "  +2.74%  "

# Row 3
$ws.Range("D3").Value = "2.457.50"
$ws.Range("E3").Value = "  +2.09%  "

# Row 4
$ws.Range("E4").Value = "  -0.12%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.80"
$ws.Range("E5").Value = "  +1.61%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.74"
$ws.Range("E6").Value = "  +2.75%  "

# Row 7
$ws.Range("E7").Value = "  +0.17%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.541"
$ws.Range("E8").Value = "  +1.04%  "

# Row 9
$ws.Range("D9").Value = "2.456.71"
$ws.Range("E9").Value = "  +1.60%  "

# Row 10
$ws.Range("E10").Value = "  +2.68%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.163"
$ws.Range("E11").Value = "  +1.82%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.27"
$ws.Range("E12").Value = "  +1.20%  "

# Row 13
$ws.Range("E13").Value = "  +2.34%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.80"
$ws.Range("E14").Value = "  +9.08%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000180"
$ws.Range("E15").Value = "  +4.06%  "

# Row 16
$ws.Range("D16").Value = "2.905.20"
$ws.Range("E16").Value = "  +2.15%  "

# Row 17
$ws.Range("D17").Value = "63.053.22"
$ws.Range("E17").Value = "  +2.92%  "

# Row 18
$ws.Range("D18").Value = "2.458.32"
$ws.Range("E18").Value = "  +1.58%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.95"
$ws.Range("E19").Value = "  -1.23%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.11"
$ws.Range("E20").Value = "  +4.43%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "329.71"
$ws.Range("E21").Value = "  +1.88%  "

# Row 22
$ws.Range("B22").Value = "SuiNetwork"
$ws.Range("C22").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.22"
$ws.Range("E22").Value = "  +12.79%  "

# Row 23
$ws.Range("B23").Value = "Polkadot"
$ws.Range("C23").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.14"
$ws.Range("E23").Value = "  +1.11%  "

# Row 24
$ws.Range("E24").Value = "  +0.06%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "66.45"
$ws.Range("E25").Value = "  +2.20%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "663.17"
$ws.Range("E26").Value = "  +8.35%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.15"
$ws.Range("E27").Value = "  +15.33%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.84"
$ws.Range("E28").Value = "  +6.75%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000100"
$ws.Range("E29").Value = "  +5.31%  "

# Row 30
$ws.Range("D30").Value = "2.576.35"

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.19"
$ws.Range("E31").Value = "  +2.32%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.43"
$ws.Range("E32").Value = "  +4.31%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.90"
$ws.Range("E33").Value = "  +5.60%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.138"
$ws.Range("E34").Value = "  +4.49%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.53"
$ws.Range("E35").Value = "  +5.03%  "

# Row 36
$ws.Range("E36").Value = "  +0.23%  "

# Row 37
$ws.Range("E37").Value = "  +3.47%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.52"
$ws.Range("E38").Value = "  +4.03%  "

# Row 39
$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "153.03"
$ws.Range("E39").Value = "  +0.20%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.373"
$ws.Range("E40").Value = "  +0.23%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.85"
$ws.Range("E41").Value = "  +2.80%  "

# Row 42
$ws.Range("B42").Value = "BabyDogeCoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D42").Value = "0.0₆0353"
$ws.Range("E42").Value = "  +23.91%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.74"
$ws.Range("E43").Value = "  +6.82%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.78"
$ws.Range("E44").Value = "  +3.92%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.31"
$ws.Range("E45").Value = "  +0.84%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.13"
$ws.Range("E47").Value = "  +28.15%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "146.55"
$ws.Range("E48").Value = "  +3.08%  "

# Row 49
$ws.Range("E49").Value = "  +2.67%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.71"
$ws.Range("E50").Value = "  +4.11%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.607"
$ws.Range("E51").Value = "  +2.37%  "
